$wb = $excel.ActiveWorkbook

# Mapping of row -> (old, new) value for column F ("想去人数") that changed
# in the "展览" sheet (sheet1) and the "全部类型" sheet (sheet4).

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 533
$ws1.Range("F6").Value  = 495
$ws1.Range("F8").Value  = 111
$ws1.Range("F10").Value = 6594
$ws1.Range("F13").Value = 2814
$ws1.Range("F14").Value = 174
$ws1.Range("F15").Value = 300
$ws1.Range("F16").Value = 255
$ws1.Range("F17").Value = 520

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 533
$ws4.Range("F8").Value  = 495
$ws4.Range("F10").Value = 111
$ws4.Range("F13").Value = 6594
$ws4.Range("F17").Value = 2814
$ws4.Range("F18").Value = 174
$ws4.Range("F19").Value = 300
$ws4.Range("F20").Value = 255
$ws4.Range("F21").Value = 520
